$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the scraped refresh
$ws.Range('D2').Value = '24.855.80'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.660.91'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'311.70"
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = "'0.3626"
$ws.Range('E7').Value = '  -1.49%  '
$ws.Range('D8').Value = "'47.26"
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('D9').Value = "'0.3261"
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('D10').Value = "'1.136"
$ws.Range('D11').Value = "'0.07058"
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').Value = "'1.001"
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = "'6.036"
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('D14').Value = "'19.53"
$ws.Range('E14').Value = '  -4.69%  '
$ws.Range('D15').Value = '1.665.49'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = "'6.610"
$ws.Range('E16').Value = '  -2.66%  '
$ws.Range('D17').Value = "'0.00001045"
$ws.Range('E17').Value = '  -4.97%  '
$ws.Range('D18').Value = "'0.06620"
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = "'79.26"
$ws.Range('E20').Value = '  -2.85%  '
$ws.Range('D21').Value = "'5.909"
$ws.Range('E21').Value = '  -4.35%  '
$ws.Range('D22').Value = "'15.72"
$ws.Range('E22').Value = '  -6.36%  '
$ws.Range('D23').Value = "'12.58"
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = '24.855.84'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').Value = "'2.423"
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').Value = "'2.408"
$ws.Range('E26').Value = '  -10.80%  '
$ws.Range('D27').Value = "'146.95"
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('D28').Value = "'18.59"
$ws.Range('E28').Value = '  -5.93%  '
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '1.849.55'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = "'1.213"
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Value = "'125.47"
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('D32').Value = "'4.117"
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('D33').Value = "'5.770"
$ws.Range('E33').Value = '  -11.27%  '
$ws.Range('D34').Value = "'0.08500"
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').Value = "'1.656"
$ws.Range('E35').Value = '  -4.07%  '
$ws.Range('D36').Value = "'12.24"
$ws.Range('E36').Value = '  -8.67%  '
$ws.Range('D37').Value = "'1.283"
$ws.Range('E37').Value = '  +4.05%  '
$ws.Range('D38').Value = "'5.151"
$ws.Range('E38').Value = '  -4.98%  '
$ws.Range('D39').Value = "'0.02250"
$ws.Range('E39').Value = '  -3.92%  '
$ws.Range('D40').Value = "'0.06029"
$ws.Range('E40').Value = '  -6.30%  '
$ws.Range('D41').Value = "'8.337"
$ws.Range('E41').Value = '  -5.41%  '
$ws.Range('D42').Value = "'0.2066"
$ws.Range('E42').Value = '  -4.15%  '
$ws.Range('D43').Value = "'0.9996"
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = "'0.5915"
$ws.Range('E44').Value = '  -5.44%  '
$ws.Range('D45').Value = "'3.776"
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = "'12.78"
$ws.Range('E46').Value = '  -4.68%  '
$ws.Range('D47').Value = "'0.5611"
$ws.Range('E47').Value = '  -5.77%  '
$ws.Range('D48').Value = "'125.02"
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').Value = "'1.944"
$ws.Range('E49').Value = '  -4.82%  '
$ws.Range('D50').Value = "'0.07000"
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').Value = "'1.194"
$ws.Range('E51').Value = '  +0.84%  '
